$d = $word.ActiveDocument

# Locate the field whose instruction text contains "m:null" (the { m:null }
# field) and remember where it lives so we can drop plain-text runs in its
# place once the field itself is gone.
$target = $null
foreach ($f in $d.Fields) {
    if ($f.Code.Text -match "null") {
        $target = $f
        break
    }
}

$para = $target.Result.Paragraphs(1)
$target.Delete()

# Replace the now-empty field paragraph with four literal text runs:
# "{", "m", ":null", "}" -- i.e. the field code spelled out as plain text
# instead of a real Word field.
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>{</w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t>:null</w:t></w:r><w:r><w:t xml:space="preserve">}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$ins = $d.Range($para.Range.Start, $para.Range.End)
[void]$ins.InsertXML($xml)
